$wb = $excel.ActiveWorkbook

# Map of worksheet name -> cell -> new header/metadata text to write.
# Every sheet's top row holds an ObjTables header string stored as a shared string;
# the sheet must be unprotected before the cell can be edited, then re-protected.
$updates = @{
    "!!Compartment" = @{
        "A1" = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='2020-03-09 23:57:49'"
        "A2" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compartment' name='Compartment' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Compound" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compound' name='Compound' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Definition" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Definition' name='Definition' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Enzyme" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Enzyme' name='Enzyme' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!FbcObjective" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='FbcObjective' name='FbcObjective' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Gene" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Gene' name='Gene' date='2020-03-09 23:57:49' objTablesVersion='0.0.8' document='lac_Operon'"
    }
    "!!Layout" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Layout' name='Layout' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Measurement" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Measurement' name='Measurement' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!PbConfig" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='PbConfig' name='PbConfig' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Position" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Position' name='Position' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Protein" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Protein' name='Protein' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Quantity" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Quantity' name='Quantity' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!QuantityInfo" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityInfo' name='QuantityInfo' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!QuantityMatrix" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityMatrix' name='QuantityMatrix' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Reaction" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Reaction' name='Reaction' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!ReactionStoichiometry" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Regulator" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Regulator' name='Regulator' date='2020-03-09 23:57:49' objTablesVersion='0.0.8' document='lac_Operon'"
    }
    "!!Relation" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relation' name='Relation' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!Relationship" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relationship' name='Relationship' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrix" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrix' name='SparseMatrix' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrixColumn" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrixOrdered" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!SparseMatrixRow" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!StoichiometricMatrix" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!rxnconContingencyList" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
    "!!rxnconReactionList" = @{
        "A1" = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconReactionList' name='rxnconReactionList' date='2020-03-09 23:57:49' objTablesVersion='0.0.8'"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Unprotect()
    $cellMap = $updates[$sheetName]
    foreach ($cellAddr in $cellMap.Keys) {
        $ws.Range($cellAddr).Value = $cellMap[$cellAddr]
    }
    $ws.Protect()
}
